# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.689.58"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "3.401.19"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'254.96"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "'654.07"
$ws.Range("E6").Value = "  +4.44%  "
$ws.Range("D7").Value = "'1.47"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.434"
$ws.Range("E8").Value = "  +4.10%  "
$ws.Range("D9").Value = "'1.07"
$ws.Range("E9").Value = "  +4.87%  "
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").Value = "3.397.80"
$ws.Range("E11").Value = "  +2.56%  "
$ws.Range("D12").Value = "'0.212"
$ws.Range("E12").Value = "  +4.09%  "
$ws.Range("D13").Value = "'41.66"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").Value = "'6.39"
$ws.Range("E14").Value = "  +18.34%  "
$ws.Range("D15").Value = "'0.0000260"
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").Value = "97.561.78"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").Value = "4.037.56"
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("D18").Value = "'8.57"
$ws.Range("E18").Value = "  +33.08%  "
$ws.Range("D19").Value = "3.408.24"
$ws.Range("E19").Value = "  +2.79%  "
$ws.Range("D20").Value = "'17.50"
$ws.Range("E20").Value = "  +10.94%  "
$ws.Range("D21").Value = "'0.491"
$ws.Range("E21").Value = "  +42.73%  "
$ws.Range("D22").Value = "'3.45"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("D23").Value = "'10.71"
$ws.Range("E23").Value = "  +13.11%  "
$ws.Range("D24").Value = "'506.89"
$ws.Range("E24").Value = "  +4.11%  "
$ws.Range("D25").Value = "'0.0000206"
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("D26").Value = "'6.18"
$ws.Range("E26").Value = "  +6.33%  "
$ws.Range("D27").Value = "'98.84"
$ws.Range("E27").Value = "  +10.59%  "
$ws.Range("D28").Value = "'12.78"
$ws.Range("E28").Value = "  +4.91%  "
$ws.Range("D29").Value = "3.592.37"
$ws.Range("E29").Value = "  +2.86%  "
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("E31").Value = "  +5.61%  "
$ws.Range("D32").Value = "'11.41"
$ws.Range("E32").Value = "  +6.63%  "
$ws.Range("D33").Value = "'0.997"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").Value = "'0.568"
$ws.Range("D36").Value = "'29.72"
$ws.Range("E36").Value = "  +5.80%  "
$ws.Range("D37").Value = "'2.29"
$ws.Range("E37").Value = "  +16.00%  "
$ws.Range("D38").Value = "'7.73"
$ws.Range("E38").Value = "  +4.87%  "
$ws.Range("D39").Value = "'528.40"
$ws.Range("E39").Value = "  +6.00%  "
$ws.Range("D40").Value = "'1.43"
$ws.Range("E40").Value = "  +13.81%  "
$ws.Range("D41").Value = "'0.153"
$ws.Range("E41").Value = "  +1.06%  "
$ws.Range("D42").Value = "'24.71"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").Value = "'0.860"
$ws.Range("E43").Value = "  +8.64%  "
$ws.Range("D44").Value = "'3.69"
$ws.Range("E44").Value = "  -6.22%  "
$ws.Range("D45").Value = "'0.0418"
$ws.Range("E45").Value = "  +21.09%  "
$ws.Range("D46").Value = "'3.31"
$ws.Range("E46").Value = "  +4.07%  "
$ws.Range("D47").Value = "'5.48"
$ws.Range("E47").Value = "  +14.90%  "
$ws.Range("D48").Value = "'8.28"
$ws.Range("E48").Value = "  +12.89%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").Value = "'1.59"
$ws.Range("E50").Value = "  +11.72%  "
$ws.Range("D51").Value = "'2.07"
$ws.Range("E51").Value = "  +5.37%  "
